$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.515.57'
$ws.Range("E2").Value = '  +3.11%  '
$ws.Range("D3").Value = '2.993.91'
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'562.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.62%  '
$ws.Range("D6").Value = "'138.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.84%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.47%  '
$ws.Range("D9").Value = '2.982.12'
$ws.Range("E9").Value = '  +2.18%  '
$ws.Range("E10").Value = '  +5.31%  '
$ws.Range("D11").Value = "'5.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +11.53%  '
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("E13").Value = '  +5.30%  '
$ws.Range("D14").Value = "'33.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.67%  '
$ws.Range("E15").Value = '  -0.22%  '
$ws.Range("D16").Value = '3.492.19'
$ws.Range("E16").Value = '  +2.55%  '
$ws.Range("E17").Value = '  +4.61%  '
$ws.Range("D18").Value = '2.990.99'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").Value = '59.522.90'
$ws.Range("E19").Value = '  +3.24%  '
$ws.Range("D20").Value = "'435.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.54%  '
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("E22").Value = '  +4.08%  '
$ws.Range("D23").Value = "'13.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.65%  '
$ws.Range("D24").Value = "'7.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.66%  '
$ws.Range("D25").Value = "'80.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.05%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  +11.22%  '
$ws.Range("E28").Value = '  +0.02%  '
$ws.Range("E29").Value = '  +3.63%  '
$ws.Range("D30").Value = "'7.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.91%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = "'6.23"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.86%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = "'25.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.47%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = "'0.106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.82%  '
$ws.Range("D34").Value = '0.0₃0777'
$ws.Range("E34").Value = '  +13.68%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = "'5.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.98%  '
$ws.Range("B36").Value = 'Mantle'
$ws.Range("C36").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D36").Value = "'0.986"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.21%  '
$ws.Range("E37").Value = '  +1.65%  '
$ws.Range("D38").Value = "'48.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.66%  '
$ws.Range("D39").Value = "'8.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.82%  '
$ws.Range("D40").Value = "'2.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.44%  '
$ws.Range("D41").Value = "'401.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.07%  '
$ws.Range("D42").Value = "'0.0354"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = '2.761.38'
$ws.Range("E43").Value = '  +2.18%  '
$ws.Range("D44").Value = "'0.105"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.59%  '
$ws.Range("D45").Value = "'0.251"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.16%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = "'123.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("D48").Value = "'34.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +19.30%  '
$ws.Range("D49").Value = "'2.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.86%  '
$ws.Range("E50").Value = '  +1.89%  '
$ws.Range("D51").Value = "'23.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.72%  '
